$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.204747333333333
$ws.Range("H2").Value = 9.614241999999999
$ws.Range("I2").Value = 0.01973032100547387
$ws.Range("J2").Value = 0.01973032100547387
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 492.021708888287
$ws.Range("R2").Value = 4428.195379994583
$ws.Range("S2").Value = 0.006258866179059004
$ws.Range("T2").Value = 0.006258866179059005
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.204747333333333
$ws.Range("H3").Value = 9.614241999999999
$ws.Range("I3").Value = 0.01973032100547387
$ws.Range("J3").Value = 0.01973032100547387
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 540.9603969884596
$ws.Range("R3").Value = 4868.643572896136
$ws.Range("S3").Value = 0.006881401108442034
$ws.Range("T3").Value = 0.006881401108442034
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.204747333333333
$ws.Range("H4").Value = 9.614241999999999
$ws.Range("I4").Value = 0.01973032100547387
$ws.Range("J4").Value = 0.01973032100547387
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 218.2122821283044
$ws.Range("R4").Value = 1963.91053915474
$ws.Range("S4").Value = 0.002775815472764476
$ws.Range("T4").Value = 0.002775815472764477
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.204747333333333
$ws.Range("H5").Value = 9.614241999999999
$ws.Range("I5").Value = 0.01973032100547387
$ws.Range("J5").Value = 0.01973032100547387
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 299.8447267962887
$ws.Range("R5").Value = 2698.602541166598
$ws.Range("S5").Value = 0.003814238245208359
$ws.Range("T5").Value = 0.003814238245208359
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 145.2141163333334
$ws.Range("H6").Value = 435.6423490000001
$ws.Range("I6").Value = 0.8940240311559332
$ws.Range("J6").Value = 0.8940240311559333
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 22294.58058358502
$ws.Range("R6").Value = 200651.2252522652
$ws.Range("S6").Value = 0.2836029261924049
$ws.Range("T6").Value = 0.283602926192405
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 145.2141163333334
$ws.Range("H7").Value = 435.6423490000001
$ws.Range("I7").Value = 0.8940240311559332
$ws.Range("J7").Value = 0.8940240311559333
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 24512.09966006942
$ws.Range("R7").Value = 220608.8969406247
$ws.Range("S7").Value = 0.3118113464683843
$ws.Range("T7").Value = 0.3118113464683843
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 145.2141163333334
$ws.Range("H8").Value = 435.6423490000001
$ws.Range("I8").Value = 0.8940240311559332
$ws.Range("J8").Value = 0.8940240311559333
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 9887.676133700948
$ws.Range("R8").Value = 88989.08520330855
$ws.Range("S8").Value = 0.1257782748703083
$ws.Range("T8").Value = 0.1257782748703083
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 145.2141163333334
$ws.Range("H9").Value = 435.6423490000001
$ws.Range("I9").Value = 0.8940240311559332
$ws.Range("J9").Value = 0.8940240311559333
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 13586.62088147963
$ws.Range("R9").Value = 122279.5879333167
$ws.Range("S9").Value = 0.1728314836248357
$ws.Range("T9").Value = 0.1728314836248357
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8052786666666667
$ws.Range("H10").Value = 2.415836
$ws.Range("I10").Value = 0.004957771998726471
$ws.Range("J10").Value = 0.004957771998726472
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 123.6336423728302
$ws.Range("R10").Value = 1112.702781355472
$ws.Range("S10").Value = 0.001572707888417328
$ws.Range("T10").Value = 0.001572707888417329
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.8052786666666667
$ws.Range("H11").Value = 2.415836
$ws.Range("I11").Value = 0.004957771998726471
$ws.Range("J11").Value = 0.004957771998726472
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 135.9307994971431
$ws.Range("R11").Value = 1223.377195474288
$ws.Range("S11").Value = 0.001729136475679952
$ws.Range("T11").Value = 0.001729136475679952
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.8052786666666667
$ws.Range("H12").Value = 2.415836
$ws.Range("I12").Value = 0.004957771998726471
$ws.Range("J12").Value = 0.004957771998726472
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 54.83168478676889
$ws.Range("R12").Value = 493.48516308092
$ws.Range("S12").Value = 0.0006974980397270468
$ws.Range("T12").Value = 0.0006974980397270471
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.8052786666666667
$ws.Range("H13").Value = 2.415836
$ws.Range("I13").Value = 0.004957771998726471
$ws.Range("J13").Value = 0.004957771998726472
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 75.34402456320935
$ws.Range("R13").Value = 678.096221068884
$ws.Range("S13").Value = 0.0009584295949021445
$ws.Range("T13").Value = 0.0009584295949021446
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.203389
$ws.Range("H14").Value = 39.610167
$ws.Range("I14").Value = 0.08128787583986632
$ws.Range("J14").Value = 0.08128787583986634
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 2027.103338639743
$ws.Range("R14").Value = 18243.93004775768
$ws.Range("S14").Value = 0.0257861966219676
$ws.Range("T14").Value = 0.02578619662196761
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.203389
$ws.Range("H15").Value = 39.610167
$ws.Range("I15").Value = 0.08128787583986632
$ws.Range("J15").Value = 0.08128787583986634
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 2228.728137392338
$ws.Range("R15").Value = 20058.55323653104
$ws.Range("S15").Value = 0.02835100750525877
$ws.Range("T15").Value = 0.02835100750525877
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.203389
$ws.Range("H16").Value = 39.610167
$ws.Range("I16").Value = 0.08128787583986632
$ws.Range("J16").Value = 0.08128787583986634
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 899.0230261057767
$ws.Range("R16").Value = 8091.207234951991
$ws.Range("S16").Value = 0.01143621248949058
$ws.Range("T16").Value = 0.01143621248949058
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 13.203389
$ws.Range("H17").Value = 39.610167
$ws.Range("I17").Value = 0.08128787583986632
$ws.Range("J17").Value = 0.08128787583986634
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 1235.344367498797
$ws.Range("R17").Value = 11118.09930748917
$ws.Range("S17").Value = 0.01571445922314937
$ws.Range("T17").Value = 0.01571445922314938
